$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-decade "artikel" counts (replacing the old algemeen/artikel/boekbespreking
# breakdown with a single, re-modeled "artikel" topic count per decade).
$counts = @(212, 211, 266, 200, 195, 204, 170, 164, 107)

# Keep only the "artikel" topic column: put the new counts into column B,
# then drop the now-unused "algemeen" and "boekbespreking" columns (C and D).
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

# Header: column B is now the "artikel" counts.
$ws.Range("B1").Value = "artikel"

# Remove the old "algemeen" (C) and "boekbespreking" (D) columns entirely.
$ws.Range("C1:D1").EntireColumn.Delete()

# Append the new decade bucket (2025, 2035] with a zero count.
$ws.Range("A11").Value = "(2025, 2035]"
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B11").Value = 0
